# TC14_Verify_AddToCart_from_Search.xlsx
# "Changes for New UI Prod"
#
# The only functional content change in this revision is on the first
# worksheet (TC14_Verify_AddToCart_from_Sear): the keyword used in row 8
# (a SCROLL_DOWN step) was replaced with a smaller-granularity scroll
# keyword, TINY_SCROLL_DOWN. Everything else in the diff (shared-string
# table physical ordering, rupBuild/revisionPtr GUIDs, the absPath of the
# author's working folder, and the customXml part numbering) is
# incidental metadata rewritten by Excel on save and is not the result of
# any user-driven edit, so it is not reproduced here.
#
# The selection (active cell) left behind on each sheet when the author
# last saved is also updated to match.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # TC14_Verify_AddToCart_from_Sear
$ws2 = $wb.Worksheets.Item(2)   # Testdata

# Update the keyword in row 8 from SCROLL_DOWN to TINY_SCROLL_DOWN
$ws1.Range("B8").Value = "TINY_SCROLL_DOWN"

# Restore the selection/active-cell that was left on each sheet
$ws2.Activate()
$ws2.Range("B2").Select()

$ws1.Activate()
$ws1.Range("B8").Select()
